$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 (Step 1 = A7, Axis 15)
$ws.Range("B7").Value = 2140
$ws.Range("D7").Value = 101.2

# Row 8 (Axis 16)
$ws.Range("B8").Value = 2140
$ws.Range("D8").Value = 96.8

# Row 9 (Axis 17)
$ws.Range("B9").Value = 2140
$ws.Range("D9").Value = 100

# Row 10 (Axis 18)
$ws.Range("B10").Value = 2140
$ws.Range("D10").Value = 95.8

$excel.Calculate()
